$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(18, 8).Value = 1331  # H18: 1085.875 -> 1331
$ws.Cells.Item(18, 9).Value = 409  # I18: 383.85715 -> 409
$ws.Cells.Item(18, 10).Value = 3175  # J18: 6000 -> 3175
$ws.Cells.Item(18, 11).Value = 409  # K18: 383.85715 -> 409
$ws.Cells.Item(18, 12).Value = 3175  # L18: 6000 -> 3175
$ws.Cells.Item(18, 13).Value = -125  # M18: -99.85714999999999 -> -125
$ws.Cells.Item(18, 14).Value = -3743  # N18: -6568 -> -3743
$ws.Cells.Item(74, 8).Value = 3992.0667  # H74: 3676.0688 -> 3992.0667
$ws.Cells.Item(74, 9).Value = 2655.9167  # I74: 2605.4614 -> 2655.9167
$ws.Cells.Item(74, 10).Value = 4882.8335  # J74: 4545.9375 -> 4882.8335
$ws.Cells.Item(74, 11).Value = 2655.9167  # K74: 2605.4614 -> 2655.9167
$ws.Cells.Item(74, 12).Value = 4882.8335  # L74: 4545.9375 -> 4882.8335
$ws.Cells.Item(74, 13).Value = -1719.9167  # M74: -1669.4614 -> -1719.9167
$ws.Cells.Item(74, 14).Value = -6754.8335  # N74: -6417.9375 -> -6754.8335
$ws.Cells.Item(77, 8).Value = 3992.0667  # H77: 3676.0688 -> 3992.0667
$ws.Cells.Item(77, 9).Value = 2655.9167  # I77: 2605.4614 -> 2655.9167
$ws.Cells.Item(77, 10).Value = 4882.8335  # J77: 4545.9375 -> 4882.8335
$ws.Cells.Item(77, 11).Value = 13279.5835  # K77: 13027.307 -> 13279.5835
$ws.Cells.Item(77, 12).Value = 24414.1675  # L77: 22729.6875 -> 24414.1675
$ws.Cells.Item(77, 13).Value = -8599.583500000001  # M77: -8347.307000000001 -> -8599.583500000001
$ws.Cells.Item(77, 14).Value = -33774.1675  # N77: -32089.6875 -> -33774.1675
$ws.Cells.Item(125, 8).Value = 5139  # H125: 5799.5713 -> 5139
$ws.Cells.Item(125, 10).Value = 5139  # J125: 5799.5713 -> 5139
$ws.Cells.Item(125, 12).Value = 46251  # L125: 52196.14169999999 -> 46251
$ws.Cells.Item(125, 14).Value = -51171  # N125: -57116.14169999999 -> -51171
$ws.Cells.Item(139, 8).Value = 147713.28  # H139: 136899 -> 147713.28
$ws.Cells.Item(139, 9).Value = 0  # I139: 130000 -> 0
$ws.Cells.Item(139, 10).Value = 147713.28  # J139: 137665.56 -> 147713.28
$ws.Cells.Item(139, 11).Value = 0  # K139: 130000 -> 0
$ws.Cells.Item(139, 12).Value = 147713.28  # L139: 137665.56 -> 147713.28
$ws.Cells.Item(139, 13).ClearContents()  # M139: -124860 -> (removed)
$ws.Cells.Item(139, 14).Value = -157993.28  # N139: -147945.56 -> -157993.28

# Sheet 2 (ARM)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 615.61536  # H2: 528.4375 -> 615.61536
$ws.Cells.Item(2, 9).Value = 545.8182  # I2: 461.14285 -> 545.8182
$ws.Cells.Item(2, 11).Value = 545.8182  # K2: 461.14285 -> 545.8182
$ws.Cells.Item(2, 13).Value = -432.8182  # M2: -348.14285 -> -432.8182
$ws.Cells.Item(32, 8).Value = 6106.147  # H32: 6285.8184 -> 6106.147
$ws.Cells.Item(32, 9).Value = 4424.759  # I32: 4576.4644 -> 4424.759
$ws.Cells.Item(32, 11).Value = 4424.759  # K32: 4576.4644 -> 4424.759
$ws.Cells.Item(32, 13).Value = -4137.759  # M32: -4289.4644 -> -4137.759
$ws.Cells.Item(43, 10).Value = 48995  # J43: 48996.332 -> 48995
$ws.Cells.Item(43, 12).Value = 48995  # L43: 48996.332 -> 48995
$ws.Cells.Item(43, 14).Value = -49621  # N43: -49622.332 -> -49621
$ws.Cells.Item(63, 8).Value = 5910.45  # H63: 6153.1055 -> 5910.45
$ws.Cells.Item(63, 9).Value = 3867.3333  # I63: 4380.8 -> 3867.3333
$ws.Cells.Item(63, 11).Value = 3867.3333  # K63: 4380.8 -> 3867.3333
$ws.Cells.Item(63, 13).Value = -3181.3333  # M63: -3694.8 -> -3181.3333
$ws.Cells.Item(66, 8).Value = 5910.45  # H66: 6153.1055 -> 5910.45
$ws.Cells.Item(66, 9).Value = 3867.3333  # I66: 4380.8 -> 3867.3333
$ws.Cells.Item(66, 11).Value = 19336.6665  # K66: 21904 -> 19336.6665
$ws.Cells.Item(66, 13).Value = -15904.6665  # M66: -18472 -> -15904.6665
$ws.Cells.Item(75, 8).Value = 70000  # H75: 69000 -> 70000
$ws.Cells.Item(75, 10).Value = 60000  # J75: 58000 -> 60000
$ws.Cells.Item(75, 12).Value = 60000  # L75: 58000 -> 60000
$ws.Cells.Item(75, 14).Value = -61748  # N75: -59748 -> -61748
$ws.Cells.Item(78, 8).Value = 70000  # H78: 69000 -> 70000
$ws.Cells.Item(78, 10).Value = 60000  # J78: 58000 -> 60000
$ws.Cells.Item(78, 12).Value = 180000  # L78: 174000 -> 180000
$ws.Cells.Item(78, 14).Value = -188736  # N78: -182736 -> -188736
$ws.Cells.Item(97, 8).Value = 2498.6667  # H97: 2498.75 -> 2498.6667
$ws.Cells.Item(97, 9).Value = 2493.9  # I97: 2548.889 -> 2493.9
$ws.Cells.Item(97, 10).Value = 2522.5  # J97: 2348.3333 -> 2522.5
$ws.Cells.Item(97, 11).Value = 2493.9  # K97: 2548.889 -> 2493.9
$ws.Cells.Item(97, 12).Value = 2522.5  # L97: 2348.3333 -> 2522.5
$ws.Cells.Item(97, 13).Value = -1997.9  # M97: -2052.889 -> -1997.9
$ws.Cells.Item(97, 14).Value = -3514.5  # N97: -3340.3333 -> -3514.5
$ws.Cells.Item(116, 8).Value = 615.61536  # H116: 528.4375 -> 615.61536
$ws.Cells.Item(116, 9).Value = 545.8182  # I116: 461.14285 -> 545.8182
$ws.Cells.Item(116, 11).Value = 545.8182  # K116: 461.14285 -> 545.8182
$ws.Cells.Item(116, 13).Value = 1748.1818  # M116: 1832.85715 -> 1748.1818

# Sheet 3 (BSM)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 615.61536  # H3: 528.4375 -> 615.61536
$ws.Cells.Item(3, 9).Value = 545.8182  # I3: 461.14285 -> 545.8182
$ws.Cells.Item(3, 11).Value = 545.8182  # K3: 461.14285 -> 545.8182
$ws.Cells.Item(3, 13).Value = -431.8182  # M3: -347.14285 -> -431.8182
$ws.Cells.Item(76, 8).Value = 120816  # H76: 152458.67 -> 120816
$ws.Cells.Item(76, 10).Value = 151088.33  # J76: 213688.5 -> 151088.33
$ws.Cells.Item(76, 12).Value = 151088.33  # L76: 213688.5 -> 151088.33
$ws.Cells.Item(76, 14).Value = -151718.33  # N76: -214318.5 -> -151718.33
$ws.Cells.Item(79, 8).Value = 120816  # H79: 152458.67 -> 120816
$ws.Cells.Item(79, 10).Value = 151088.33  # J79: 213688.5 -> 151088.33
$ws.Cells.Item(79, 12).Value = 151088.33  # L79: 213688.5 -> 151088.33
$ws.Cells.Item(79, 14).Value = -153272.33  # N79: -215872.5 -> -153272.33
$ws.Cells.Item(86, 8).Value = 18100.611  # H86: 18111.334 -> 18100.611
$ws.Cells.Item(86, 9).Value = 9601.666999999999  # I86: 9617.75 -> 9601.666999999999
$ws.Cells.Item(86, 11).Value = 9601.666999999999  # K86: 9617.75 -> 9601.666999999999
$ws.Cells.Item(86, 13).Value = -8478.666999999999  # M86: -8494.75 -> -8478.666999999999
$ws.Cells.Item(89, 8).Value = 18100.611  # H89: 18111.334 -> 18100.611
$ws.Cells.Item(89, 9).Value = 9601.666999999999  # I89: 9617.75 -> 9601.666999999999
$ws.Cells.Item(89, 11).Value = 48008.335  # K89: 48088.75 -> 48008.335
$ws.Cells.Item(89, 13).Value = -42392.335  # M89: -42472.75 -> -42392.335
$ws.Cells.Item(94, 8).Value = 1181.8667  # H94: 1218.5 -> 1181.8667
$ws.Cells.Item(94, 9).Value = 798.3  # I94: 812.6667 -> 798.3
$ws.Cells.Item(94, 11).Value = 798.3  # K94: 812.6667 -> 798.3
$ws.Cells.Item(94, 13).Value = -347.3  # M94: -361.6667 -> -347.3
$ws.Cells.Item(122, 8).Value = 84780  # H122: 0 -> 84780
$ws.Cells.Item(122, 10).Value = 84780  # J122: 0 -> 84780
$ws.Cells.Item(122, 12).Value = 84780  # L122: 0 -> 84780
$ws.Cells.Item(122, 14).Value = -94580  # N122: None -> -94580
$ws.Cells.Item(130, 8).Value = 50000  # H130: 49996.668 -> 50000
$ws.Cells.Item(130, 10).Value = 50000  # J130: 49996.668 -> 50000
$ws.Cells.Item(130, 12).Value = 50000  # L130: 49996.668 -> 50000
$ws.Cells.Item(130, 14).Value = -60040  # N130: -60036.668 -> -60040
$ws.Cells.Item(134, 8).Value = 3905.8  # H134: 4439 -> 3905.8
$ws.Cells.Item(134, 9).Value = 3827.6428  # I134: 4392.25 -> 3827.6428
$ws.Cells.Item(134, 11).Value = 11482.9284  # K134: 13176.75 -> 11482.9284
$ws.Cells.Item(134, 13).Value = -8947.928400000001  # M134: -10641.75 -> -8947.928400000001

# Sheet 4 (CRP)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(86, 8).Value = 4514.7  # H86: 4422.4546 -> 4514.7
$ws.Cells.Item(86, 9).Value = 4250  # I86: 4125 -> 4250
$ws.Cells.Item(86, 11).Value = 4250  # K86: 4125 -> 4250
$ws.Cells.Item(86, 13).Value = -3127  # M86: -3002 -> -3127
$ws.Cells.Item(89, 8).Value = 4514.7  # H89: 4422.4546 -> 4514.7
$ws.Cells.Item(89, 9).Value = 4250  # I89: 4125 -> 4250
$ws.Cells.Item(89, 11).Value = 21250  # K89: 20625 -> 21250
$ws.Cells.Item(89, 13).Value = -15634  # M89: -15009 -> -15634
$ws.Cells.Item(107, 8).Value = 2158  # H107: 2255.6924 -> 2158
$ws.Cells.Item(107, 9).Value = 1051.8  # I107: 1070 -> 1051.8
$ws.Cells.Item(107, 11).Value = 1051.8  # K107: 1070 -> 1051.8
$ws.Cells.Item(107, 13).Value = 868.2  # M107: 850 -> 868.2
$ws.Cells.Item(132, 8).Value = 4801.268  # H132: 4735.756 -> 4801.268
$ws.Cells.Item(132, 10).Value = 8999.75  # J132: 8328.25 -> 8999.75
$ws.Cells.Item(132, 12).Value = 26999.25  # L132: 24984.75 -> 26999.25
$ws.Cells.Item(132, 14).Value = -32059.25  # N132: -30044.75 -> -32059.25

# Sheet 5 (CUL)
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(22, 8).Value = 725.9375  # H22: 651.2222 -> 725.9375
$ws.Cells.Item(22, 9).Value = 141  # I22: 135.75 -> 141
$ws.Cells.Item(22, 10).Value = 9500  # J22: 4775 -> 9500
$ws.Cells.Item(22, 11).Value = 423  # K22: 407.25 -> 423
$ws.Cells.Item(22, 12).Value = 28500  # L22: 14325 -> 28500
$ws.Cells.Item(22, 13).Value = -254  # M22: -238.25 -> -254
$ws.Cells.Item(22, 14).Value = -28838  # N22: -14663 -> -28838
$ws.Cells.Item(27, 8).Value = 725.9375  # H27: 651.2222 -> 725.9375
$ws.Cells.Item(27, 9).Value = 141  # I27: 135.75 -> 141
$ws.Cells.Item(27, 10).Value = 9500  # J27: 4775 -> 9500
$ws.Cells.Item(27, 11).Value = 423  # K27: 407.25 -> 423
$ws.Cells.Item(27, 12).Value = 28500  # L27: 14325 -> 28500
$ws.Cells.Item(27, 13).Value = -321  # M27: -305.25 -> -321
$ws.Cells.Item(27, 14).Value = -28704  # N27: -14529 -> -28704
$ws.Cells.Item(44, 8).Value = 3066.125  # H44: 3081.111 -> 3066.125
$ws.Cells.Item(44, 9).Value = 2066  # I44: 2218.5715 -> 2066
$ws.Cells.Item(44, 10).Value = 4733  # J44: 6100 -> 4733
$ws.Cells.Item(44, 11).Value = 6198  # K44: 6655.7145 -> 6198
$ws.Cells.Item(44, 12).Value = 14199  # L44: 18300 -> 14199
$ws.Cells.Item(44, 13).Value = -5800  # M44: -6257.7145 -> -5800
$ws.Cells.Item(44, 14).Value = -14995  # N44: -19096 -> -14995
$ws.Cells.Item(47, 8).Value = 1780.6  # H47: 1475.75 -> 1780.6
$ws.Cells.Item(47, 10).Value = 2666.6667  # J47: 2500 -> 2666.6667
$ws.Cells.Item(47, 12).Value = 8000.000100000001  # L47: 7500 -> 8000.000100000001
$ws.Cells.Item(47, 14).Value = -8862.000100000001  # N47: -8362 -> -8862.000100000001
$ws.Cells.Item(98, 8).Value = 1887.8889  # H98: 2156.375 -> 1887.8889
$ws.Cells.Item(98, 9).Value = 916  # I98: 1151.2 -> 916
$ws.Cells.Item(98, 11).Value = 2748  # K98: 3453.6 -> 2748
$ws.Cells.Item(98, 13).Value = -1250  # M98: -1955.6 -> -1250
$ws.Cells.Item(132, 8).Value = 1636.7778  # H132: 1583.6364 -> 1636.7778
$ws.Cells.Item(132, 10).Value = 1585.5  # J132: 1505.1666 -> 1585.5
$ws.Cells.Item(132, 12).Value = 14269.5  # L132: 13546.4994 -> 14269.5
$ws.Cells.Item(132, 14).Value = -19329.5  # N132: -18606.4994 -> -19329.5
$ws.Cells.Item(141, 8).Value = 8021.44  # H141: 7980.231 -> 8021.44
$ws.Cells.Item(141, 10).Value = 11553.2  # J141: 11134.728 -> 11553.2
$ws.Cells.Item(141, 12).Value = 34659.60000000001  # L141: 33404.18399999999 -> 34659.60000000001
$ws.Cells.Item(141, 14).Value = -45019.60000000001  # N141: -43764.18399999999 -> -45019.60000000001

# Sheet 6 (GSM)
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(33, 8).Value = 9997  # H33: 9996 -> 9997
$ws.Cells.Item(33, 10).Value = 9997  # J33: 9996 -> 9997
$ws.Cells.Item(33, 12).Value = 9997  # L33: 9996 -> 9997
$ws.Cells.Item(33, 14).Value = -10501  # N33: -10500 -> -10501
$ws.Cells.Item(80, 8).Value = 4698.778  # H80: 4722.222 -> 4698.778
$ws.Cells.Item(80, 9).Value = 4612.857  # I80: 4915.8335 -> 4612.857
$ws.Cells.Item(80, 10).Value = 4999.5  # J80: 4335 -> 4999.5
$ws.Cells.Item(80, 11).Value = 4612.857  # K80: 4915.8335 -> 4612.857
$ws.Cells.Item(80, 12).Value = 4999.5  # L80: 4335 -> 4999.5
$ws.Cells.Item(80, 13).Value = -3614.857  # M80: -3917.8335 -> -3614.857
$ws.Cells.Item(80, 14).Value = -6995.5  # N80: -6331 -> -6995.5
$ws.Cells.Item(83, 8).Value = 4698.778  # H83: 4722.222 -> 4698.778
$ws.Cells.Item(83, 9).Value = 4612.857  # I83: 4915.8335 -> 4612.857
$ws.Cells.Item(83, 10).Value = 4999.5  # J83: 4335 -> 4999.5
$ws.Cells.Item(83, 11).Value = 23064.285  # K83: 24579.1675 -> 23064.285
$ws.Cells.Item(83, 12).Value = 24997.5  # L83: 21675 -> 24997.5
$ws.Cells.Item(83, 13).Value = -18072.285  # M83: -19587.1675 -> -18072.285
$ws.Cells.Item(83, 14).Value = -34981.5  # N83: -31659 -> -34981.5
$ws.Cells.Item(102, 8).Value = 2688.087  # H102: 2427.1155 -> 2688.087
$ws.Cells.Item(102, 9).Value = 2025.5883  # I102: 1785.7 -> 2025.5883
$ws.Cells.Item(102, 11).Value = 2025.5883  # K102: 1785.7 -> 2025.5883
$ws.Cells.Item(102, 13).Value = -403.5882999999999  # M102: -163.7 -> -403.5882999999999
$ws.Cells.Item(132, 8).Value = 3023.5881  # H132: 2506.0466 -> 3023.5881
$ws.Cells.Item(132, 9).Value = 2818.8262  # I132: 2419.1428 -> 2818.8262
$ws.Cells.Item(132, 10).Value = 3451.7273  # J132: 2668.2666 -> 3451.7273
$ws.Cells.Item(132, 11).Value = 8456.4786  # K132: 7257.428400000001 -> 8456.4786
$ws.Cells.Item(132, 12).Value = 10355.1819  # L132: 8004.7998 -> 10355.1819
$ws.Cells.Item(132, 13).Value = -5926.4786  # M132: -4727.428400000001 -> -5926.4786
$ws.Cells.Item(132, 14).Value = -15415.1819  # N132: -13064.7998 -> -15415.1819

# Sheet 7 (LTW)
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(23, 8).Value = 499  # H23: 0 -> 499
$ws.Cells.Item(23, 9).Value = 499  # I23: 0 -> 499
$ws.Cells.Item(23, 11).Value = 499  # K23: 0 -> 499
$ws.Cells.Item(23, 13).Value = -269  # M23: None -> -269
$ws.Cells.Item(33, 8).Value = 0  # H33: 15000 -> 0
$ws.Cells.Item(33, 10).Value = 0  # J33: 15000 -> 0
$ws.Cells.Item(33, 12).Value = 0  # L33: 15000 -> 0
$ws.Cells.Item(33, 14).ClearContents()  # N33: -15580 -> (removed)
$ws.Cells.Item(93, 8).Value = 2497.0454  # H93: 2583.0952 -> 2497.0454
$ws.Cells.Item(93, 9).Value = 1807.625  # I93: 1967.2858 -> 1807.625
$ws.Cells.Item(93, 11).Value = 1807.625  # K93: 1967.2858 -> 1807.625
$ws.Cells.Item(93, 13).Value = -559.625  # M93: -719.2858000000001 -> -559.625
$ws.Cells.Item(132, 8).Value = 0  # H132: 399.5 -> 0
$ws.Cells.Item(132, 9).Value = 0  # I132: 399.5 -> 0
$ws.Cells.Item(132, 11).Value = 0  # K132: 1198.5 -> 0
$ws.Cells.Item(132, 13).ClearContents()  # M132: 1331.5 -> (removed)
$ws.Cells.Item(136, 8).Value = 5312.0625  # H136: 4888.5 -> 5312.0625
$ws.Cells.Item(136, 9).Value = 2000  # I136: 1875 -> 2000
$ws.Cells.Item(136, 11).Value = 6000  # K136: 5625 -> 6000
$ws.Cells.Item(136, 13).Value = -3450  # M136: -3075 -> -3450

# Sheet 8 (WVR)
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(45, 8).Value = 29984.25  # H45: 29995 -> 29984.25
$ws.Cells.Item(45, 10).Value = 29984.25  # J45: 29995 -> 29984.25
$ws.Cells.Item(45, 12).Value = 29984.25  # L45: 29995 -> 29984.25
$ws.Cells.Item(45, 14).Value = -30966.25  # N45: -30977 -> -30966.25
$ws.Cells.Item(68, 8).Value = 52667  # H68: 54999.668 -> 52667
$ws.Cells.Item(68, 9).Value = 54000  # I68: 52500 -> 54000
$ws.Cells.Item(68, 10).Value = 50001  # J68: 59999 -> 50001
$ws.Cells.Item(68, 11).Value = 54000  # K68: 52500 -> 54000
$ws.Cells.Item(68, 12).Value = 50001  # L68: 59999 -> 50001
$ws.Cells.Item(68, 13).Value = -53189  # M68: -51689 -> -53189
$ws.Cells.Item(68, 14).Value = -51623  # N68: -61621 -> -51623
$ws.Cells.Item(69, 8).Value = 0  # H69: 30000 -> 0
$ws.Cells.Item(69, 10).Value = 0  # J69: 30000 -> 0
$ws.Cells.Item(69, 12).Value = 0  # L69: 30000 -> 0
$ws.Cells.Item(69, 14).ClearContents()  # N69: -31498 -> (removed)
$ws.Cells.Item(71, 8).Value = 52667  # H71: 54999.668 -> 52667
$ws.Cells.Item(71, 9).Value = 54000  # I71: 52500 -> 54000
$ws.Cells.Item(71, 10).Value = 50001  # J71: 59999 -> 50001
$ws.Cells.Item(71, 11).Value = 162000  # K71: 157500 -> 162000
$ws.Cells.Item(71, 12).Value = 150003  # L71: 179997 -> 150003
$ws.Cells.Item(71, 13).Value = -157944  # M71: -153444 -> -157944
$ws.Cells.Item(71, 14).Value = -158115  # N71: -188109 -> -158115
$ws.Cells.Item(72, 8).Value = 0  # H72: 30000 -> 0
$ws.Cells.Item(72, 10).Value = 0  # J72: 30000 -> 0
$ws.Cells.Item(72, 12).Value = 0  # L72: 90000 -> 0
$ws.Cells.Item(72, 14).ClearContents()  # N72: -97488 -> (removed)
$ws.Cells.Item(110, 8).Value = 0  # H110: 84995 -> 0
$ws.Cells.Item(110, 10).Value = 0  # J110: 84995 -> 0
$ws.Cells.Item(110, 12).Value = 0  # L110: 84995 -> 0
$ws.Cells.Item(110, 14).ClearContents()  # N110: -93175 -> (removed)
